# Triathlon Season workbook update
#
# 1) Renames league names to their full "... League" form.
# 2) Fixes a handful of event/partner-product description typos and
#    wording tweaks across the existing leagues.
# 3) Fills in the previously-blank Hunter League rows 38-40 (only the
#    Clubs column had data) and appends seven brand-new rows to finish
#    the Hunter League, plus three new leagues: Central West League,
#    New England League and Hume League (rows 41-56).
#
# Values that are purely numeric-looking text (the "Round" column) are
# written with a leading apostrophe so Excel stores them as text, matching
# the rest of the sheet instead of coercing them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- South Coast League (rows 2-6): "south coast" -> "South Coast League"
$ws.Range("A2").Value = "South Coast League"
$ws.Range("A3").Value = "South Coast League"
$ws.Range("A4").Value = "South Coast League"
$ws.Range("F4").Value = "Sprint Aquabike, Aquabike, Super Sprint"
$ws.Range("A5").Value = "South Coast League"
$ws.Range("E5").Value = "Super Sprint, Sprint, Classic and Ironman 70.3"
$ws.Range("A6").Value = "South Coast League"

# --- North Coast League (rows 7-13): "North Coast" -> "North Coast League"
$ws.Range("A7").Value = "North Coast League"
$ws.Range("A8").Value = "North Coast League"
$ws.Range("F8").Value = "Super Sprint, Aquathon, Teams"
$ws.Range("A9").Value = "North Coast League"
$ws.Range("A10").Value = "North Coast League"
$ws.Range("A11").Value = "North Coast League"
$ws.Range("F11").Value = "Sprint Aquabike, Aquabike, Super Sprint"
$ws.Range("A12").Value = "North Coast League"
$ws.Range("A13").Value = "North Coast League"

# --- Sydney Premier League wording fixes (rows 14, 17, 18, 19)
$ws.Range("E14").Value = "Ironman 70.3 and Sprint"
$ws.Range("F14").Value = "n/a"
$ws.Range("E17").Value = "Long Aqua"
$ws.Range("F17").Value = "Short Aqua"
$ws.Range("E18").Value = "Super Sprint, Sprint"
$ws.Range("E19").Value = "Super Sprint, Sprint, Classic and Ironman 70.3"

# --- Sydney League 1 wording fixes (rows 24, 25, 27, 28, 29, 30)
$ws.Range("E24").Value = "Ironman 70.3 and Sprint"
$ws.Range("F24").Value = "n/a"
$ws.Range("E25").Value = "Standard, Aquabike"
$ws.Range("E27").Value = "Long Aqua"
$ws.Range("F27").Value = "Short Aqua"
$ws.Range("E28").Value = "Super Sprint, Sprint"
$ws.Range("E29").Value = "Super Sprint, Sprint, Classic and Ironman 70.3"
$ws.Range("E30").Value = "Sprint, Standard"

# --- Hunter League existing rows 36-37 corrections
$ws.Range("C36").Value = "Sparke Helmore Triathlon"
$ws.Range("E36").Value = "Sprint and Standard"
$ws.Range("C37").Value = "Hawks Nest Triathlon"
$ws.Range("E37").Value = "Sprint, Standard"
$ws.Range("F37").Value = "Super Sprint, Aquabike"

# --- Hunter League rows 38-40: previously only had the Clubs (G) column
$ws.Range("A38").Value = "Hunter League"
$ws.Range("B38").Value = "'3"
$ws.Range("C38").Value = "Singleton Triathlon"
$ws.Range("D38").Value = "No"
$ws.Range("E38").Value = "Sprint"
$ws.Range("F38").Value = "Super Sprint"

$ws.Range("A39").Value = "Hunter League"
$ws.Range("B39").Value = "'4"
$ws.Range("C39").Value = "CCTRI Club Race"
$ws.Range("D39").Value = "No"
$ws.Range("E39").Value = "Sprint"
$ws.Range("F39").Value = "Super Sprint"

$ws.Range("A40").Value = "Hunter League"
$ws.Range("B40").Value = "'5"
$ws.Range("C40").Value = "Island Triathlon"
$ws.Range("D40").Value = "No"
$ws.Range("E40").Value = "Super Sprint, Sprint, Standard, Aquabike"
$ws.Range("F40").Value = "n/a"

# --- New Hunter League rows 41-42 (no Clubs entry for either)
$ws.Range("A41").Value = "Hunter League"
$ws.Range("B41").Value = "'6"
$ws.Range("C41").Value = "NSW Triathlon Club Champs (Double Points)"
$ws.Range("D41").Value = "yes"
$ws.Range("E41").Value = "Sprint, Standard "
$ws.Range("F41").Value = "Sprint Aquabike, Aquabike, Super Sprint"

$ws.Range("A42").Value = "Hunter League"
$ws.Range("B42").Value = "'7"
$ws.Range("C42").Value = "Lake Macquarie Triathlon Festival"
$ws.Range("D42").Value = "No"
$ws.Range("E42").Value = "Sprint and Standard"
$ws.Range("F42").Value = "Super Sprint"

# --- New Central West League rows 43-47
$ws.Range("A43").Value = "Central West League"
$ws.Range("B43").Value = "'1"
$ws.Range("C43").Value = "Mudgee Triathlon Club"
$ws.Range("D43").Value = "No"
$ws.Range("E43").Value = "Club Distance"
$ws.Range("F43").Value = "n/a"
$ws.Range("G43").Value = "Mudgee Triathlon Club"

$ws.Range("A44").Value = "Central West League"
$ws.Range("B44").Value = "'2"
$ws.Range("C44").Value = "Bathurst  Triathlon Club"
$ws.Range("D44").Value = "No"
$ws.Range("E44").Value = "Club Distance"
$ws.Range("F44").Value = "n/a"
$ws.Range("G44").Value = "Bathurst  Triathlon Club"

$ws.Range("A45").Value = "Central West League"
$ws.Range("B45").Value = "'3"
$ws.Range("C45").Value = "Orange Triathlon Club"
$ws.Range("D45").Value = "No"
$ws.Range("E45").Value = "Club Distance"
$ws.Range("F45").Value = "n/a"
$ws.Range("G45").Value = "Orange Triathlon Club"

$ws.Range("A46").Value = "Central West League"
$ws.Range("B46").Value = "'4"
$ws.Range("C46").Value = "Dubbo Triathlon Club"
$ws.Range("D46").Value = "No"
$ws.Range("E46").Value = "Club Distance"
$ws.Range("F46").Value = "n/a"
$ws.Range("G46").Value = "Dubbo Triathlon Club"

# Row 47 only has a Clubs entry, no round data
$ws.Range("G47").Value = "Cowra Triathlon Club"

# --- New England League rows 48-51
$ws.Range("A48").Value = "New England League"
$ws.Range("B48").Value = "'1"
$ws.Range("C48").Value = "Scone Club Race "
$ws.Range("D48").Value = "No"
$ws.Range("E48").Value = "Sprint and Super Sprint"
$ws.Range("F48").Value = "n/a"
$ws.Range("G48").Value = "Scone Triathlon Club"

$ws.Range("A49").Value = "New England League"
$ws.Range("B49").Value = "'2"
$ws.Range("C49").Value = "Tamworth Club Race"
$ws.Range("D49").Value = "No"
$ws.Range("E49").Value = "Club Distance"
$ws.Range("F49").Value = "n/a"
$ws.Range("G49").Value = "Armidale Triathlon Club"

$ws.Range("A50").Value = "New England League"
$ws.Range("B50").Value = "'3"
$ws.Range("C50").Value = "Armidale Club Race"
$ws.Range("D50").Value = "No"
$ws.Range("E50").Value = "Club Distance"
$ws.Range("F50").Value = "n/a"
$ws.Range("G50").Value = "Gunnedah Triathlon Club"

$ws.Range("A51").Value = "New England League"
$ws.Range("B51").Value = "'4"
$ws.Range("C51").Value = "NSW Triathlon Club Champs"
$ws.Range("D51").Value = "yes"
$ws.Range("E51").Value = "Sprint, Standard"
$ws.Range("F51").Value = "Sprint Aquabike, Aquabike, Super Sprint"
$ws.Range("G51").Value = "Tamworth Triathlon Club"

# --- New Hume League rows 52-56
$ws.Range("A52").Value = "Hume League"
$ws.Range("B52").Value = "'1"
$ws.Range("C52").Value = "Challenge Canberra"
$ws.Range("D52").Value = "No"
$ws.Range("E52").Value = "Ironman 70.3, Aquabike 70.3, Standard, Aquabike, Sprint"
$ws.Range("F52").Value = "Super Sprint"
$ws.Range("G52").Value = "Yass Valley Triathlon Club"

$ws.Range("A53").Value = "Hume League"
$ws.Range("B53").Value = "'2"
$ws.Range("C53").Value = "Jackie Fairweather Triathlon"
$ws.Range("D53").Value = "No"
$ws.Range("E53").Value = "Sprint"
$ws.Range("F53").Value = "Super Sprint"
$ws.Range("G53").Value = "Goulburn Triathlon Club"

$ws.Range("A54").Value = "Hume League"
$ws.Range("B54").Value = "'3"
$ws.Range("C54").Value = "MMJ Aquathlon"
$ws.Range("D54").Value = "No"
$ws.Range("E54").Value = "Long Aqua"
$ws.Range("F54").Value = "Short Aqua, Mini Aqua"
$ws.Range("G54").Value = "Leeton Triathlon Club"

$ws.Range("A55").Value = "Hume League"
$ws.Range("B55").Value = "'4"
$ws.Range("C55").Value = "NSW Triathlon Club Champs (Double Points)"
$ws.Range("D55").Value = "yes"
$ws.Range("E55").Value = "Sprint, Standard"
$ws.Range("F55").Value = "Sprint Aquabike, Aquabike, Super Sprint"
$ws.Range("G55").Value = "Wagga TrIathlon Club"

$ws.Range("A56").Value = "Hume League"
$ws.Range("B56").Value = "'5"
$ws.Range("C56").Value = "Big Husky"
$ws.Range("D56").Value = "No"
$ws.Range("E56").Value = "Super Sprint, Sprint, Classic and Ironman 70.3"
$ws.Range("F56").Value = "Aquabike"
